# Regenerate save_data to use K instead of Strike#, updating the K column
# (column G) values on Sheet1 for each trade row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "K" values keyed by row number (row 1 is the header row).
$kValues = @{
    2  = 1
    4  = 2
    5  = 2
    7  = 0
    8  = 3
    9  = 1
    10 = 1
    11 = 1
    12 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
